# Search Course.xlsx - mark search test cases TC1:01..TC1:09 as passing.
# For each of those rows the "actual result" (column E) now matches the
# "expected result" (column D), and the Result column (F) flips from
# "Failed" to "Pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Range("D$r").Value = $ws.Range("E$r").Value()
    $ws.Range("F$r").Value = "Pass"
}

# Row 3's ExpectedResult cell (D3) had lost its cell border in the old
# data; restore the same formatting used by the rest of the column by
# copying D2's format onto it.
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
[void]$ws.Application.CutCopyMode

# Reflect the reviewer's final selection (they had just finished marking
# the Result column as Pass for rows 2-10).
[void]$ws.Range("F2:F10").Select()
